$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.002.36'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -4.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.742.70'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.84%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.61'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5805'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2733'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.25'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06623'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07554'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.743.46'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.711'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6037'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.981.23'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '74.67'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008729'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -11.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.005.89'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.332'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.81%  '
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '205.62'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.30'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.635'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.08%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.138'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1235'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.381'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06159'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.04%  '
$ws.Range("E31").Value = '  -3.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.746'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.729'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.670'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.038'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6410'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.417'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.721'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01669'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.132.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.187'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8771'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.892.10'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.49'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.582'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("E48").Value = '  -5.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.298'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05385'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.269'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.65%  '
